$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: a new week's observation is inserted at row 228 (just
# below the most-recent row 227), and every subsequent historical row
# shifts down by one to make room, with the former last row (247)
# becoming the new last row (248).

function Copy-Row($src, $dst) {
    for ($col = 1; $col -le 18; $col++) {
        $srcCell = $ws.Cells.Item($src, $col)
        $dstCell = $ws.Cells.Item($dst, $col)
        $dstCell.Value = $srcCell.Value2()
    }
    # Column D carries the date serial; keep its date number format.
    $ws.Cells.Item($dst, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# Shift rows down by one, starting from the bottom so we never clobber
# a row before it has been copied.
for ($r = 247; $r -ge 228; $r--) {
    $dst = $r + 1
    Copy-Row $r $dst
}

# New data for row 228.
$ws.Cells.Item(228, 1).Value = 11
$ws.Cells.Item(228, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(228, 3).Value = "Bíobío"
$ws.Cells.Item(228, 4).Value = 45223
$ws.Cells.Item(228, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(228, 5).Value = 8
$ws.Cells.Item(228, 6).Value = 100112043
$ws.Cells.Item(228, 7).Value = "Pepino ensalada"
$ws.Cells.Item(228, 8).Value = "Sin especificar"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 120
$ws.Cells.Item(228, 11).Value = 14000
$ws.Cells.Item(228, 12).Value = 14000
$ws.Cells.Item(228, 13).Value = 14000
$ws.Cells.Item(228, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(228, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(228, 16).Value = 233
$ws.Cells.Item(228, 17).Value = 60
$ws.Cells.Item(228, 18).Value = "Hortaliza"
